$d = $word.ActiveDocument
$d.Content.Find.Execute("Edward Venator", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Edward Venator TEST", 2)
